# "Add files via upload" -- adds a new "Fixation Data" worksheet (with its
# summary data) to the workbook, as the third/last sheet, and leaves it as
# the active sheet/tab when the file is saved -- matching the state the
# workbook was left in after the upload.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Scroll/select on "Task Answers" the way it was left before switching
#    to the new sheet (selection ends up on I2, columns scrolled so C is
#    the left-most visible column).
# ---------------------------------------------------------------------
$wsTaskAnswers = $wb.Worksheets.Item("Task Answers")
$wsTaskAnswers.Activate()
$null = $wsTaskAnswers.Range("I2").Select()

# ---------------------------------------------------------------------
# 2) Add the new "Fixation Data" worksheet after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFixation = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsFixation.Name = "Fixation Data"

# Header row
$wsFixation.Range("B1").Value = "Task 1"
$wsFixation.Range("C1").Value = "Task 2"
$wsFixation.Range("D1").Value = "Task 3"
$wsFixation.Range("E1").Value = "Task 4"

# Fixations Per Person
$wsFixation.Range("A2").Value = "Fixations Per Person"
$wsFixation.Range("B2").Value = 326
$wsFixation.Range("C2").Value = "err"
$wsFixation.Range("D2").Value = "err"
$wsFixation.Range("E2").Value = 680

# Time Per Person (duration values, displayed h:mm)
$wsFixation.Range("A3").Value = "Time Per Person"
$wsFixation.Range("B3").Value = 0.12638888888888888
$wsFixation.Range("C3").Value = 0.19097222222222221
$wsFixation.Range("D3").Value = 0.26111111111111113
$wsFixation.Range("E3").Value = 0.25416666666666665
$wsFixation.Range("B3:E3").NumberFormat = "h:mm"

# Average Fixations Per Second
$wsFixation.Range("A4").Value = "Average Fixations Per Second"
$wsFixation.Range("B4").Value = 1.79
$wsFixation.Range("C4").Value = "err"
$wsFixation.Range("D4").Value = "err"
$wsFixation.Range("E4").Value = 1.86

# ---------------------------------------------------------------------
# 3) Leave the new sheet active, selected at E6 (matches the saved file).
# ---------------------------------------------------------------------
$wsFixation.Activate()
$null = $wsFixation.Range("E6").Select()
